# Atualizacao de bases das ligas, do dia: 19-06-2024 as 21:51
# Two pairs of rows had their match data (columns B through AD) swapped
# between each other, while the running index in column A stayed put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    for ($col = 2; $col -le 30; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)

        $val1 = $cell1.Value()
        $val2 = $cell2.Value()

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

# Rows 49 and 50 (A=47 and A=48)
Swap-RowData 49 50

# Rows 190 and 191 (A=188 and A=189)
Swap-RowData 190 191
